# Regenerate the per-trial stimulus list + sanity-check columns for the
# "bedrooms" memory block (subject 12): trial_total renumbered, and the
# target/new/catch rows + their conceptual/perceptual/typicality/n/p_*
# columns (H:S) are rewritten to the finalized input list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: trial_total 55 -> 27; stimulus -> 'stimuli/img_ic3os.png'
$ws.Cells.Item(2, 6).Value = 27
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_ic3os.png"
$row[0,5] = 84.79069767441861
$row[0,6] = 66.16279069767442
$row[0,7] = 75.47674418604652
$row[0,8] = 43
$row[0,9] = 9
$row[0,10] = 9
$row[0,11] = 9
$ws.Range("H2:S2").Value = $row

# Row 3: trial_total 56 -> 28; stimulus -> 'stimuli/img_jivhq.png'
$ws.Cells.Item(3, 6).Value = 28
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_jivhq.png"
$row[0,5] = 37
$row[0,6] = 22.26530612244898
$row[0,7] = 29.63265306122449
$row[0,8] = 49
$row[0,9] = 2
$row[0,10] = 2
$row[0,11] = 2
$ws.Range("H3:S3").Value = $row

# Row 4: trial_total 57 -> 29; stimulus -> 'stimuli/img_1vq1v.png'
$ws.Cells.Item(4, 6).Value = 29
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_1vq1v.png"
$row[0,5] = 69.42857142857143
$row[0,6] = 46.59523809523809
$row[0,7] = 58.01190476190476
$row[0,8] = 42
$row[0,9] = 5
$row[0,10] = 5
$row[0,11] = 5
$ws.Range("H4:S4").Value = $row

# Row 5: trial_total 58 -> 30; stimulus -> 'stimuli/img_72fmj.png'
$ws.Cells.Item(5, 6).Value = 30
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_72fmj.png"
$row[0,5] = 53.87179487179487
$row[0,6] = 36.02564102564103
$row[0,7] = 44.94871794871795
$row[0,8] = 39
$row[0,9] = 3
$row[0,10] = 3
$row[0,11] = 3
$ws.Range("H5:S5").Value = $row

# Row 6: trial_total 59 -> 31; stimulus -> 'stimuli/img_t4hvr.png'
$ws.Cells.Item(6, 6).Value = 31
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_t4hvr.png"
$row[0,5] = 61.69230769230769
$row[0,6] = 39.76923076923077
$row[0,7] = 50.73076923076923
$row[0,8] = 39
$row[0,9] = 3
$row[0,10] = 3
$row[0,11] = 3
$ws.Range("H6:S6").Value = $row

# Row 7: trial_total 60 -> 32; stimulus -> 'stimuli/img_f4jxo.png'
$ws.Cells.Item(7, 6).Value = 32
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_f4jxo.png"
$row[0,5] = 82.91666666666667
$row[0,6] = 65.52777777777777
$row[0,7] = 74.22222222222223
$row[0,8] = 36
$row[0,9] = 8
$row[0,10] = 8
$row[0,11] = 8
$ws.Range("H7:S7").Value = $row

# Row 8: trial_total 61 -> 33; stimulus -> 'stimuli/img_sltwe.png'
$ws.Cells.Item(8, 6).Value = 33
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_sltwe.png"
$row[0,5] = 72.025
$row[0,6] = 46.875
$row[0,7] = 59.45
$row[0,8] = 40
$row[0,9] = 5
$row[0,10] = 5
$row[0,11] = 5
$ws.Range("H8:S8").Value = $row

# Row 9: trial_total 62 -> 34; stimulus -> 'stimuli/img_c4uwt.png'
$ws.Cells.Item(9, 6).Value = 34
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_c4uwt.png"
$row[0,5] = 44.48387096774194
$row[0,6] = 30.06451612903226
$row[0,7] = 37.2741935483871
$row[0,8] = 31
$row[0,9] = 2
$row[0,10] = 2
$row[0,11] = 2
$ws.Range("H9:S9").Value = $row

# Row 10: trial_total 63 -> 35; stimulus -> 'stimuli/catch_08.jpg'
$ws.Cells.Item(10, 6).Value = 35
$row = New-Object 'object[,]' 1,12
$row[0,0] = $null
$row[0,1] = $null
$row[0,2] = "catch"
$row[0,3] = "f"
$row[0,4] = "stimuli/catch_08.jpg"
$row[0,5] = $null
$row[0,6] = $null
$row[0,7] = $null
$row[0,8] = $null
$row[0,9] = $null
$row[0,10] = $null
$row[0,11] = $null
$ws.Range("H10:S10").Value = $row

# Row 11: trial_total 64 -> 36; stimulus -> 'stimuli/img_ose78.png'
$ws.Cells.Item(11, 6).Value = 36
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_ose78.png"
$row[0,5] = 80.19444444444444
$row[0,6] = 60.25
$row[0,7] = 70.22222222222223
$row[0,8] = 36
$row[0,9] = 8
$row[0,10] = 7
$row[0,11] = 7
$ws.Range("H11:S11").Value = $row

# Row 12: trial_total 65 -> 37; stimulus -> 'stimuli/img_cgdyc.png'
$ws.Cells.Item(12, 6).Value = 37
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_cgdyc.png"
$row[0,5] = 32.93023255813954
$row[0,6] = 14.04651162790698
$row[0,7] = 23.48837209302326
$row[0,8] = 43
$row[0,9] = 1
$row[0,10] = 1
$row[0,11] = 1
$ws.Range("H12:S12").Value = $row

# Row 13: trial_total 66 -> 38; stimulus -> 'stimuli/img_aweye.png'
$ws.Cells.Item(13, 6).Value = 38
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_aweye.png"
$row[0,5] = 53.42105263157895
$row[0,6] = 31.84210526315789
$row[0,7] = 42.63157894736842
$row[0,8] = 38
$row[0,9] = 2
$row[0,10] = 2
$row[0,11] = 2
$ws.Range("H13:S13").Value = $row

# Row 14: trial_total 67 -> 39; stimulus -> 'stimuli/img_zi682.png'
$ws.Cells.Item(14, 6).Value = 39
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_zi682.png"
$row[0,5] = 84.6
$row[0,6] = 69.525
$row[0,7] = 77.0625
$row[0,8] = 40
$row[0,9] = 9
$row[0,10] = 9
$row[0,11] = 9
$ws.Range("H14:S14").Value = $row

# Row 15: trial_total 68 -> 40; stimulus -> 'stimuli/img_jge7p.png'
$ws.Cells.Item(15, 6).Value = 40
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_jge7p.png"
$row[0,5] = 90.42424242424242
$row[0,6] = 75.63636363636364
$row[0,7] = 83.03030303030303
$row[0,8] = 33
$row[0,9] = 10
$row[0,10] = 10
$row[0,11] = 10
$ws.Range("H15:S15").Value = $row

# Row 16: trial_total 69 -> 41; stimulus -> 'stimuli/img_2pnl2.png'
$ws.Cells.Item(16, 6).Value = 41
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_2pnl2.png"
$row[0,5] = 6.621621621621622
$row[0,6] = 7.135135135135135
$row[0,7] = 6.878378378378379
$row[0,8] = 37
$row[0,9] = 1
$row[0,10] = 1
$row[0,11] = 1
$ws.Range("H16:S16").Value = $row

# Row 17: trial_total 70 -> 42; stimulus -> 'stimuli/img_yteqw.png'
$ws.Cells.Item(17, 6).Value = 42
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_yteqw.png"
$row[0,5] = 66.83783783783784
$row[0,6] = 43.78378378378378
$row[0,7] = 55.31081081081081
$row[0,8] = 37
$row[0,9] = 4
$row[0,10] = 4
$row[0,11] = 4
$ws.Range("H17:S17").Value = $row

# Row 18: trial_total 71 -> 43; stimulus -> 'stimuli/img_kzg3h.png'
$ws.Cells.Item(18, 6).Value = 43
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_kzg3h.png"
$row[0,5] = 77.02777777777777
$row[0,6] = 56.22222222222222
$row[0,7] = 66.625
$row[0,8] = 36
$row[0,9] = 7
$row[0,10] = 7
$row[0,11] = 7
$ws.Range("H18:S18").Value = $row

# Row 19: trial_total 72 -> 44; stimulus -> 'stimuli/img_anzgh.png'
$ws.Cells.Item(19, 6).Value = 44
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_anzgh.png"
$row[0,5] = 75.10526315789474
$row[0,6] = 55.76315789473684
$row[0,7] = 65.4342105263158
$row[0,8] = 38
$row[0,9] = 6
$row[0,10] = 6
$row[0,11] = 6
$ws.Range("H19:S19").Value = $row

# Row 20: trial_total 73 -> 45; stimulus -> 'stimuli/img_cmyvx.png'
$ws.Cells.Item(20, 6).Value = 45
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_cmyvx.png"
$row[0,5] = 64.25
$row[0,6] = 40.09375
$row[0,7] = 52.171875
$row[0,8] = 32
$row[0,9] = 4
$row[0,10] = 4
$row[0,11] = 4
$ws.Range("H20:S20").Value = $row

# Row 21: trial_total 74 -> 46; stimulus -> 'stimuli/img_fnu4h.png'
$ws.Cells.Item(21, 6).Value = 46
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_fnu4h.png"
$row[0,5] = 85.87179487179488
$row[0,6] = 70.71794871794872
$row[0,7] = 78.2948717948718
$row[0,8] = 39
$row[0,9] = 9
$row[0,10] = 9
$row[0,11] = 9
$ws.Range("H21:S21").Value = $row

# Row 22: trial_total 75 -> 47; stimulus -> 'stimuli/img_9pfbj.png'
$ws.Cells.Item(22, 6).Value = 47
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_9pfbj.png"
$row[0,5] = 91.27272727272727
$row[0,6] = 80.0909090909091
$row[0,7] = 85.68181818181819
$row[0,8] = 33
$row[0,9] = 10
$row[0,10] = 10
$row[0,11] = 10
$ws.Range("H22:S22").Value = $row

# Row 23: trial_total 76 -> 48; stimulus -> 'stimuli/img_z3yzz.png'
$ws.Cells.Item(23, 6).Value = 48
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_z3yzz.png"
$row[0,5] = 71.71052631578948
$row[0,6] = 49.81578947368421
$row[0,7] = 60.76315789473685
$row[0,8] = 38
$row[0,9] = 5
$row[0,10] = 5
$row[0,11] = 5
$ws.Range("H23:S23").Value = $row

# Row 24: trial_total 77 -> 49; stimulus -> 'stimuli/img_gbypq.png'
$ws.Cells.Item(24, 6).Value = 49
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_gbypq.png"
$row[0,5] = 76.275
$row[0,6] = 51.925
$row[0,7] = 64.1
$row[0,8] = 40
$row[0,9] = 6
$row[0,10] = 6
$row[0,11] = 6
$ws.Range("H24:S24").Value = $row

# Row 25: trial_total 78 -> 50; stimulus -> 'stimuli/img_ozxpp.png'
$ws.Cells.Item(25, 6).Value = 50
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_ozxpp.png"
$row[0,5] = 26.26470588235294
$row[0,6] = 11.47058823529412
$row[0,7] = 18.86764705882353
$row[0,8] = 34
$row[0,9] = 1
$row[0,10] = 1
$row[0,11] = 1
$ws.Range("H25:S25").Value = $row

# Row 26: trial_total 79 -> 51; stimulus -> 'stimuli/img_3bxjb.png'
$ws.Cells.Item(26, 6).Value = 51
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_3bxjb.png"
$row[0,5] = 87.28571428571429
$row[0,6] = 72.65714285714286
$row[0,7] = 79.97142857142858
$row[0,8] = 35
$row[0,9] = 10
$row[0,10] = 10
$row[0,11] = 10
$ws.Range("H26:S26").Value = $row

# Row 27: trial_total 80 -> 52; stimulus -> 'stimuli/img_psgf7.png'
$ws.Cells.Item(27, 6).Value = 52
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_psgf7.png"
$row[0,5] = 26
$row[0,6] = 11.66666666666667
$row[0,7] = 18.83333333333333
$row[0,8] = 36
$row[0,9] = 1
$row[0,10] = 1
$row[0,11] = 1
$ws.Range("H27:S27").Value = $row

# Row 28: trial_total 81 -> 53; stimulus -> 'stimuli/img_4wq98.png'
$ws.Cells.Item(28, 6).Value = 53
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = $null
$row[0,2] = "new"
$row[0,3] = "f"
$row[0,4] = "stimuli/img_4wq98.png"
$row[0,5] = 78.48387096774194
$row[0,6] = 58.12903225806452
$row[0,7] = 68.30645161290323
$row[0,8] = 31
$row[0,9] = 7
$row[0,10] = 7
$row[0,11] = 7
$ws.Range("H28:S28").Value = $row

# Row 29: trial_total 82 -> 54; stimulus -> 'stimuli/img_juob3.png'
$ws.Cells.Item(29, 6).Value = 54
$row = New-Object 'object[,]' 1,12
$row[0,0] = "bedrooms"
$row[0,1] = "target"
$row[0,2] = "old"
$row[0,3] = "j"
$row[0,4] = "stimuli/img_juob3.png"
$row[0,5] = 79.92105263157895
$row[0,6] = 59.78947368421053
$row[0,7] = 69.85526315789474
$row[0,8] = 38
$row[0,9] = 7
$row[0,10] = 7
$row[0,11] = 7
$ws.Range("H29:S29").Value = $row

